# Auto-generated script to apply scheduled runner price/profit updates
# to the Famfrit_Profits workbook (8 job sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

### Sheet: ALC ###
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 2926.8948
$ws.Range("I19").Value = 1398.4445
$ws.Range("K19").Value = 1398.4445
$ws.Range("M19").Value = -1223.4445
# Row 32
$ws.Range("H32").Value = 5077.125
$ws.Range("I32").Value = 4108
$ws.Range("J32").Value = 6046.25
$ws.Range("K32").Value = 4108
$ws.Range("L32").Value = 6046.25
$ws.Range("M32").Value = -3782
$ws.Range("N32").Value = -6698.25
# Row 33
$ws.Range("H33").Value = 299.5
$ws.Range("I33").Value = 308.45456
$ws.Range("J33").Value = 201
$ws.Range("K33").Value = 308.45456
$ws.Range("L33").Value = 201
$ws.Range("M33").Value = -79.45456000000001
$ws.Range("N33").Value = -659
# Row 48
$ws.Range("H48").Value = 3416.6667
$ws.Range("I48").Value = 3416.6667
$ws.Range("K48").Value = 10250.0001
$ws.Range("M48").Value = -9958.000100000001
# Row 56
$ws.Range("H56").Value = 3416.6667
$ws.Range("I56").Value = 3416.6667
$ws.Range("K56").Value = 10250.0001
$ws.Range("M56").Value = -9716.000100000001
# Row 70
$ws.Range("H70").Value = 2255.2222
$ws.Range("I70").Value = 1499
$ws.Range("J70").Value = 2471.2856
$ws.Range("K70").Value = 4497
$ws.Range("L70").Value = 7413.8568
$ws.Range("M70").Value = -4227
$ws.Range("N70").Value = -7953.8568
# Row 73
$ws.Range("H73").Value = 2255.2222
$ws.Range("I73").Value = 1499
$ws.Range("J73").Value = 2471.2856
$ws.Range("K73").Value = 4497
$ws.Range("L73").Value = 7413.8568
$ws.Range("M73").Value = -3561
$ws.Range("N73").Value = -9285.856800000001
# Row 107
$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080
# Row 109
$ws.Range("H109").Value = 45100
$ws.Range("J109").Value = 40200
$ws.Range("L109").Value = 40200
$ws.Range("N109").Value = -42974
# Row 112
$ws.Range("H112").Value = 27781218
$ws.Range("J112").Value = 27781218
$ws.Range("L112").Value = 83343654
$ws.Range("N112").Value = -83345870
# Row 121
$ws.Range("H121").Value = 1294.75
$ws.Range("J121").Value = 1294.75
$ws.Range("L121").Value = 3884.25
$ws.Range("N121").Value = -7378.25
# Row 132
$ws.Range("H132").Value = 2636.2058
$ws.Range("I132").Value = 2216.3215
$ws.Range("K132").Value = 6648.9645
$ws.Range("M132").Value = -4118.9645
# Row 137
$ws.Range("H137").Value = 3021.0527
$ws.Range("I137").Value = 3112.7646
$ws.Range("K137").Value = 9338.293799999999
$ws.Range("M137").Value = -6788.293799999999

### Sheet: ARM ###
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 333365980
$ws.Range("J43").Value = 48997.5
$ws.Range("L43").Value = 48997.5
$ws.Range("N43").Value = -49623.5
# Row 74
$ws.Range("H74").Value = 125002760
$ws.Range("I74").Value = 200001600
$ws.Range("J74").Value = 4689.3335
$ws.Range("K74").Value = 200001600
$ws.Range("L74").Value = 4689.3335
$ws.Range("M74").Value = -200000726
$ws.Range("N74").Value = -6437.3335
# Row 77
$ws.Range("H77").Value = 125002760
$ws.Range("I77").Value = 200001600
$ws.Range("J77").Value = 4689.3335
$ws.Range("K77").Value = 1000008000
$ws.Range("L77").Value = 23446.6675
$ws.Range("M77").Value = -1000003632
$ws.Range("N77").Value = -32182.6675
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

### Sheet: BSM ###
$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

### Sheet: CRP ###
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 3965.3333
$ws.Range("I107").Value = 2001
$ws.Range("J107").Value = 4947.5
$ws.Range("K107").Value = 2001
$ws.Range("L107").Value = 4947.5
$ws.Range("M107").Value = -81
$ws.Range("N107").Value = -8787.5

### Sheet: CUL ###
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 28337382
$ws.Range("J4").Value = 5308333.5
$ws.Range("L4").Value = 15925000.5
$ws.Range("N4").Value = -15925224.5
# Row 80
$ws.Range("H80").Value = 4999.7617
$ws.Range("I80").Value = 4998.3335
$ws.Range("K80").Value = 14995.0005
$ws.Range("M80").Value = -14059.0005
# Row 83
$ws.Range("H83").Value = 4999.7617
$ws.Range("I83").Value = 4998.3335
$ws.Range("K83").Value = 44985.0015
$ws.Range("M83").Value = -40305.0015
# Row 99
$ws.Range("H99").Value = 756.3333
$ws.Range("I99").Value = 756.3333
$ws.Range("K99").Value = 2268.9999
$ws.Range("M99").Value = -22.9998999999998
# Row 121
$ws.Range("H121").Value = 1213.6364
$ws.Range("I121").Value = 150
$ws.Range("K121").Value = 450
$ws.Range("M121").Value = 860
# Row 137
$ws.Range("H137").Value = 5527.1113
$ws.Range("I137").Value = 3589.4
$ws.Range("J137").Value = 6272.385
$ws.Range("K137").Value = 10768.2
$ws.Range("L137").Value = 18817.155
$ws.Range("M137").Value = -5668.200000000001
$ws.Range("N137").Value = -29017.155
# Row 138
$ws.Range("H138").Value = 3419
$ws.Range("J138").Value = 3688.6667
$ws.Range("L138").Value = 11066.0001
$ws.Range("N138").Value = -21346.0001

### Sheet: LTW ###
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 698.9
$ws.Range("I16").Value = 643.2222
$ws.Range("K16").Value = 643.2222
$ws.Range("M16").Value = -473.2222
# Row 46
$ws.Range("H46").Value = 1239.6182
$ws.Range("I46").Value = 514.4358999999999
$ws.Range("J46").Value = 3007.25
$ws.Range("K46").Value = 514.4358999999999
$ws.Range("L46").Value = 3007.25
$ws.Range("M46").Value = -326.4358999999999
$ws.Range("N46").Value = -3383.25
# Row 93
$ws.Range("H93").Value = 2056.3333
$ws.Range("I93").Value = 806.4706
$ws.Range("K93").Value = 806.4706
$ws.Range("M93").Value = 441.5294
# Row 136
$ws.Range("H136").Value = 2245.068
$ws.Range("I136").Value = 1726.5366
$ws.Range("J136").Value = 9331.666999999999
$ws.Range("K136").Value = 5179.6098
$ws.Range("L136").Value = 27995.001
$ws.Range("M136").Value = -2629.6098
$ws.Range("N136").Value = -33095.001

### Sheet: WVR ###
$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 2999
$ws.Range("I11").Value = 2999
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 2999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -2857
$ws.Range("N11").ClearContents()
# Row 136
$ws.Range("H136").Value = 1976.8
$ws.Range("I136").Value = 1769.2106
$ws.Range("K136").Value = 5307.6318
$ws.Range("M136").Value = -2757.6318

Write-Host "Applied scheduled price/profit updates to all sheets."